$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.400.57'
$ws.Range("E2").Value = '  -0.99%  '
$ws.Range("D3").Value = '1.571.24'
$ws.Range("E3").Value = '  -0.80%  '
$ws.Range("E4").Value = '  -0.12%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '206.97'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.498'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.87%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '22.15'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.249'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.80%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0594'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.71%  '
$ws.Range("E11").Value = '  +0.03%  '
$ws.Range("D12").Value = '1.797.13'
$ws.Range("E12").Value = '  -0.62%  '
$ws.Range("D13").Value = '1.569.74'
$ws.Range("E13").Value = '  -0.34%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.81'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.518'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '63.45'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.01%  '
$ws.Range("D17").Value = '27.415.81'
$ws.Range("E17").Value = '  -0.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '213.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -2.77%  '
$ws.Range("D19").Value = '0.0₃0691'
$ws.Range("E19").Value = '  -0.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.27'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.53%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.14'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '9.48'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("E24").Value = '  +1.35%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '153.18'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.43%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '6.67'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -2.58%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '14.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -1.07%  '
$ws.Range("E28").Value = '  -0.86%  '
$ws.Range("E29").Value = '  -0.13%  '
$ws.Range("E30").Value = '  -0.44%  '
$ws.Range("E31").Value = '  +1.27%  '
$ws.Range("E32").Value = '  -0.37%  '
$ws.Range("D33").Value = '1.393.78'
$ws.Range("E33").Value = '  +2.47%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.97'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +1.10%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.55'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.75%  '
$ws.Range("E36").Value = '  -0.46%  '
$ws.Range("E37").Value = '  -2.78%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0166'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.37%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.529'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.92%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.821'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.33%  '
$ws.Range("E41").Value = '  -0.11%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.995'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.44%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.81'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '64.10'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.69%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.18'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.38%  '
$ws.Range("E46").Value = '  +1.47%  '
$ws.Range("D47").Value = '1.708.35'
$ws.Range("E47").Value = '  -0.67%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '85.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -2.60%  '
$ws.Range("D49").Value = '0.0₇0996'
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0494'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.62%  '
$ws.Range("B51").Value = 'Algorand'
$ws.Range("C51").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0951'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.83%  '
